$d = $word.ActiveDocument

# Replace the curly-apostrophe contraction "That’s" with "That is"
# (the run containing "\u2019" becomes " i", merging with neighboring
# text so the paragraph reads "That is all for this file.")
$d.Content.Find.Execute([ref]"hat’s", [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]"hat is", [ref]2)
